$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1599003333333333
$ws.Range("H2").Value = 0.479701
$ws.Range("I2").Value = 0.0264777194346773
$ws.Range("J2").Value = 0.02647771943467731
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 31.38723566666667
$ws.Range("N2").Value = 94.161707
$ws.Range("O2").Value = 0.5539598599114094
$ws.Range("P2").Value = 0.5539598599114095
$ws.Range("Q2").Value = 5.01882944551189
$ws.Range("R2").Value = 45.169465009607
$ws.Range("S2").Value = 0.01466759374880744
$ws.Range("T2").Value = 0.01466759374880745

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1599003333333333
$ws.Range("H3").Value = 0.479701
$ws.Range("I3").Value = 0.0264777194346773
$ws.Range("J3").Value = 0.02647771943467731
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.55525033333333
$ws.Range("N3").Value = 52.665751
$ws.Range("O3").Value = 0.3098362697066353
$ws.Range("P3").Value = 0.3098362697066353
$ws.Range("Q3").Value = 2.807090380050111
$ws.Range("R3").Value = 25.263813420451
$ws.Range("S3").Value = 0.008203757819979297
$ws.Range("T3").Value = 0.008203757819979297

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1599003333333333
$ws.Range("H4").Value = 0.479701
$ws.Range("I4").Value = 0.0264777194346773
$ws.Range("J4").Value = 0.02647771943467731
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.717279333333334
$ws.Range("N4").Value = 23.151838
$ws.Range("O4").Value = 0.1362038703819552
$ws.Range("P4").Value = 0.1362038703819552
$ws.Range("Q4").Value = 1.233995537826445
$ws.Range("R4").Value = 11.105959840438
$ws.Range("S4").Value = 0.003606367865890564
$ws.Range("T4").Value = 0.003606367865890565

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.879152
$ws.Range("H5").Value = 17.637456
$ws.Range("I5").Value = 0.9735222805653226
$ws.Range("J5").Value = 0.9735222805653228
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 31.38723566666667
$ws.Range("N5").Value = 94.161707
$ws.Range("O5").Value = 0.5539598599114094
$ws.Range("P5").Value = 0.5539598599114095
$ws.Range("Q5").Value = 184.5303293441547
$ws.Range("R5").Value = 1660.772964097392
$ws.Range("S5").Value = 0.539292266162602
$ws.Range("T5").Value = 0.5392922661626021

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.879152
$ws.Range("H6").Value = 17.637456
$ws.Range("I6").Value = 0.9735222805653226
$ws.Range("J6").Value = 0.9735222805653228
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.55525033333333
$ws.Range("N6").Value = 52.665751
$ws.Range("O6").Value = 0.3098362697066353
$ws.Range("P6").Value = 0.3098362697066353
$ws.Range("Q6").Value = 103.2099851077173
$ws.Range("R6").Value = 928.889865969456
$ws.Range("S6").Value = 0.301632511886656
$ws.Range("T6").Value = 0.301632511886656

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.879152
$ws.Range("H7").Value = 17.637456
$ws.Range("I7").Value = 0.9735222805653226
$ws.Range("J7").Value = 0.9735222805653228
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.717279333333334
$ws.Range("N7").Value = 23.151838
$ws.Range("O7").Value = 0.1362038703819552
$ws.Range("P7").Value = 0.1362038703819552
$ws.Range("Q7").Value = 45.37105822712534
$ws.Range("R7").Value = 408.339524044128
$ws.Range("S7").Value = 0.1325975025160646
$ws.Range("T7").Value = 0.1325975025160647

